$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 2
$ws.Cells.Item(2, 5).Value = 94

# Row 10
$ws.Cells.Item(10, 5).Value = 518
$ws.Cells.Item(10, 6).Value = 257
$ws.Cells.Item(10, 7).Value = 94
$ws.Cells.Item(10, 8).Value = 351

# Row 11
$ws.Cells.Item(11, 5).Value = 337
$ws.Cells.Item(11, 7).Value = 68
$ws.Cells.Item(11, 8).Value = 254

# Row 12
$ws.Cells.Item(12, 5).Value = 506
$ws.Cells.Item(12, 6).Value = 276
$ws.Cells.Item(12, 8).Value = 360

# Row 14
$ws.Cells.Item(14, 5).Value = 125

# Row 15
$ws.Cells.Item(15, 5).Value = 160

# Row 16
$ws.Cells.Item(16, 5).Value = 198

# Row 21
$ws.Cells.Item(21, 5).Value = 137

# Row 23
$ws.Cells.Item(23, 5).Value = 198

# Row 24
$ws.Cells.Item(24, 5).Value = 208
$ws.Cells.Item(24, 6).Value = 113
$ws.Cells.Item(24, 8).Value = 143

# Row 25
$ws.Cells.Item(25, 5).Value = 258
$ws.Cells.Item(25, 6).Value = 126
$ws.Cells.Item(25, 8).Value = 186

# Row 27
$ws.Cells.Item(27, 5).Value = 320
$ws.Cells.Item(27, 6).Value = 160
$ws.Cells.Item(27, 8).Value = 241

# Row 28
$ws.Cells.Item(28, 5).Value = 193
$ws.Cells.Item(28, 6).Value = 75
$ws.Cells.Item(28, 8).Value = 127

# Row 30
$ws.Cells.Item(30, 5).Value = 201

# Row 33
$ws.Cells.Item(33, 5).Value = 287
$ws.Cells.Item(33, 6).Value = 145
$ws.Cells.Item(33, 8).Value = 234

# Row 37
$ws.Cells.Item(37, 5).Value = 154

# Row 42
$ws.Cells.Item(42, 5).Value = 372
$ws.Cells.Item(42, 6).Value = 205
$ws.Cells.Item(42, 8).Value = 265

# Row 44
$ws.Cells.Item(44, 5).Value = 309

# Row 45
$ws.Cells.Item(45, 5).Value = 140
$ws.Cells.Item(45, 6).Value = 69
$ws.Cells.Item(45, 8).Value = 108

# Row 46
$ws.Cells.Item(46, 5).Value = 311
$ws.Cells.Item(46, 6).Value = 169
$ws.Cells.Item(46, 8).Value = 232

# Row 47
$ws.Cells.Item(47, 5).Value = 443
$ws.Cells.Item(47, 6).Value = 223
$ws.Cells.Item(47, 8).Value = 315

# Row 48
$ws.Cells.Item(48, 5).Value = 202
$ws.Cells.Item(48, 6).Value = 88
$ws.Cells.Item(48, 8).Value = 132

# Row 49
$ws.Cells.Item(49, 5).Value = 281

# Row 50
$ws.Cells.Item(50, 5).Value = 238
$ws.Cells.Item(50, 6).Value = 113
$ws.Cells.Item(50, 8).Value = 186

# Row 51
$ws.Cells.Item(51, 5).Value = 227
$ws.Cells.Item(51, 6).Value = 99
$ws.Cells.Item(51, 8).Value = 171

# Row 52
$ws.Cells.Item(52, 5).Value = 26

$wb.Save()
